$wb = $excel.ActiveWorkbook

# --- Update the existing "TradeInDevice" sheet ---
$trade = $wb.Worksheets.Item("TradeInDevice")

# A2 changes from "1" to "0"
$trade.Range("A2").Value = "0"

# Column F (6) gets an explicit (custom) width, same visual width as the default
$trade.Columns.Item(6).ColumnWidth = 10

# Selection moves to D2
$trade.Range("D2").Select() | Out-Null

# --- Add the new "AppleTablets" sheet after the last sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$apple = $wb.Worksheets.Add($null, $lastSheet)
$apple.Name = "AppleTablets"
$apple.Range("A2").Value = "Shop Apple tablets"

# Reuse the header text style already used elsewhere in the workbook
$styled = $wb.Worksheets.Item("Google5GPhones")
$styled.Range("A2").Copy() | Out-Null
$apple.Range("A2").PasteSpecial(-4122) | Out-Null

# Column A width
$apple.Columns.Item(1).ColumnWidth = 21

# Selection on the new sheet
$apple.Range("A2").Select() | Out-Null
